# "Generate Report for Archive"
#
# The localization-status report is regenerated: entries that were
# previously marked "Ready for handoff" have moved on to "In Translation".
# Because the new status text is shorter than the old one, the "Status"
# columns (zh-cn/de-de values on the Overview sheet, and the Status column
# on each per-locale sheet) shrink to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears ("Ready for handoff" -> "In Translation")
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# Re-fit the affected "Status" columns now that their content is shorter
$overview.Range("E:F").EntireColumn.ColumnWidth = 12.5
$zhcn.Range("C:C").EntireColumn.ColumnWidth = 12.5
$dede.Range("C:C").EntireColumn.ColumnWidth = 12.5
